$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "311.69"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.95%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.29"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.15%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.143"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.07%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07819"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.38%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.901"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.20%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.265"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.81%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.836"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-9.24%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9193"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.22%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1193"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.81%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1922"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.25%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08994"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.01%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03342"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.00%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09593"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.07%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001383"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.05%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005705"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.68%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.532"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.52%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.418"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.44%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.00%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.269"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "4.89%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1284"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.65%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04361"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.69%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001253"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.11%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004667"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.93%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001365"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.81%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004000"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02290"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.48%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05051"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.90%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007467"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.17%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009079"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-8.42%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.20%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.001957"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.28%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009417"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "10.59%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006690"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.19%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.01%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.001003"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-23.07%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003288"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "9.53%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002106"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.01%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.01%"
